$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "Arrete l'incendie" (fire quest) updated to 2-maisons variant ---
$ws.Range("F4").Value2 = "1`n5`n4`n2`n2`n3`n2`n1"
$ws.Range("G4").Value2 = "2`n4`n6`n2`n1`n3`n2`n1"
$ws.Range("H4").Value2 = "1`n1`n1`n2`n2`n3`n2`n1"
$ws.Range("I4").Value2 = "Il faut apporter 4000mL d'eau pour la maison 1`nIl faut apporter 1200mL d'eau pour la maison 2`nIl faut apporter 5200mL d'eau au total"
$ws.Range("J4").Value2 = "Il faut apporter 9600mL d'eau pour la maison 1`nIl faut apporter 600mL d'eau pour la maison 2`nIl faut apporter 10200mL d'eau au total"
$ws.Range("K4").Value2 = "Il faut apporter 200mL d'eau pour la maison 1`nIl faut apporter 1200mL d'eau pour la maison 2`nIl faut apporter 1400mL d'eau au total"
$ws.Range("L4").Value2 = "Vite tu dois arreter l'incendie`nTu vas recevoir les dimensions de 2 maisons par 2 maisons`nNavi est partie te dire l'étendue des flammes pendant que tu es parti cherché de l'eau`nNavi va te dire les informations dans l'ordre suivant:`nla puissance des flammes`nla longueur de la maison`nla largeur de la maison`nla nombre d'étage de la maison`nElle répetera ces informations pour les 2 maisons`nAprès chaque maison, tu devras lui dire quel quantité d'eau ammener`nUne fois les 2 maisons éteintes, tu devras annoncer quel quantité d'eau vous avez utilisés`nUne maison avec un incendie de puissance 1, de longueur 1, de largeur 1 et de 1 étage demandera 100mL d'eau "
$ws.Range("R4").Value2 = "1`n5`n4`n2`n2`n3`n2`n1"
$ws.Range("S4").Value2 = "Il faut apporter 4000mL d'eau pour la maison 1`nIl faut apporter 1200mL d'eau pour la maison 2`nIl faut apporter 5200mL d'eau au total"
$ws.Range("T4").Value2 = "5`n2`n1`n5`n3`n1`n6`n8"
$ws.Range("U4").Value2 = "Il faut apporter 5000mL d'eau pour la maison 1`nIl faut apporter 14400mL d'eau pour la maison 2`nIl faut apporter 19400mL d'eau au total"
$ws.Range("V4").Value2 = "2`n1`n6`n1`n4`n3`n1`n2"
$ws.Range("W4").Value2 = "Il faut apporter 1200mL d'eau pour la maison 1`nIl faut apporter 2400mL d'eau pour la maison 2`nIl faut apporter 3600mL d'eau au total"

# --- Row 5: "Dechiffrer le message" -> "Trouve ton chemin" (new riddle quest) ---
$ws.Range("B5").Value2 = "Trouve ton chemin"
$ws.Range("F5").Value2 = 60
$ws.Range("G5").Value2 = 41
$ws.Range("H5").Value2 = 84
$ws.Range("I5").Value2 = "60 est divisible par`n2`n3"
$ws.Range("J5").Value2 = "41 n'est divisible par aucun des trois nombres"
$ws.Range("K5").Value2 = "84 est divisible par`n2`n3`n7"
$ws.Range("L5").Value2 = "Devant toi se trouve trois portes`nEn reseolvant cette enigme, vous trouvez la porte à prendre`nPour un entier n donné`nDeterminer si celui çi est pair,s'il est divisible par 3 et ou s'il est divisible par 7 ou aucun des 3`nExemple : Pour un entier 42,`nLa réponse sera : `n42 est divisible par`n2`n3`n7"
$ws.Range("Q5").Value2 = 4
$ws.Range("R5").Value2 = 60
$ws.Range("S5").Value2 = "60 est divisible par`n2`n3"
$ws.Range("T5").Value2 = 43
$ws.Range("U5").Value2 = "43 n'est divisible par aucun des trois nombres"
$ws.Range("V5").Value2 = 84
$ws.Range("W5").Value2 = "40 est divisible par`n2"
$ws.Range("X5").Value2 = 168
$ws.Range("Y5").Value2 = "168 est divisible par`n2`n3`n7"
